# Scheduled market-data refresh: rewrites the computed crafting-profit
# columns (currentAveragePrice / NQ / HQ / LevePrice NQ+HQ / LeveProfit
# NQ+HQ, columns H-N) on the affected rows of each Leve sheet with the
# latest pulled marketboard figures.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2300.6667
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 902
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 902
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -1252
# Row 58
$ws.Range("H58").Value = 1902
$ws.Range("J58").Value = 2285.5715
$ws.Range("L58").Value = 6856.7145
$ws.Range("N58").Value = -7156.7145
# Row 76
$ws.Range("H76").Value = 6944.75
$ws.Range("I76").Value = 7008.231
$ws.Range("K76").Value = 7008.231
$ws.Range("M76").Value = -6693.231
# Row 79
$ws.Range("H79").Value = 6944.75
$ws.Range("I79").Value = 7008.231
$ws.Range("K79").Value = 7008.231
$ws.Range("M79").Value = -5916.231
# Row 80
$ws.Range("H80").Value = 1672653.8
$ws.Range("J80").Value = 2564345.5
$ws.Range("L80").Value = 7693036.5
$ws.Range("N80").Value = -7695032.5
# Row 83
$ws.Range("H83").Value = 1672653.8
$ws.Range("J83").Value = 2564345.5
$ws.Range("L83").Value = 23079109.5
$ws.Range("N83").Value = -23089093.5
# Row 116
$ws.Range("H116").Value = 3009.4814
$ws.Range("I116").Value = 2652.96
$ws.Range("K116").Value = 2652.96
$ws.Range("M116").Value = 789.04
# Row 125
$ws.Range("H125").Value = 8765.556
$ws.Range("J125").Value = 13333
$ws.Range("L125").Value = 119997
$ws.Range("N125").Value = -124917
# Row 138
$ws.Range("H138").Value = 3655.913
$ws.Range("I138").Value = 4932.846
$ws.Range("J138").Value = 3152.879
$ws.Range("K138").Value = 14798.538
$ws.Range("L138").Value = 9458.636999999999
$ws.Range("M138").Value = -9658.537999999999
$ws.Range("N138").Value = -19738.637

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5275625
$ws.Range("I32").Value = 7250697
$ws.Range("K32").Value = 7250697
$ws.Range("M32").Value = -7250410
# Row 43
$ws.Range("H43").Value = 61286.332
$ws.Range("J43").Value = 61286.332
$ws.Range("L43").Value = 61286.332
$ws.Range("N43").Value = -61912.332
# Row 45
$ws.Range("H45").Value = 4294.9653
$ws.Range("I45").Value = 3982.8845
$ws.Range("J45").Value = 6999.6665
$ws.Range("K45").Value = 3982.8845
$ws.Range("L45").Value = 6999.6665
$ws.Range("M45").Value = -3605.8845
$ws.Range("N45").Value = -7753.6665
# Row 92
$ws.Range("H92").Value = 90550
$ws.Range("J92").Value = 90550
$ws.Range("L92").Value = 90550
$ws.Range("N92").Value = -95542
# Row 122
$ws.Range("H122").Value = 3665.2068
$ws.Range("I122").Value = 2406.9375
$ws.Range("K122").Value = 7220.8125
$ws.Range("M122").Value = -4770.8125
# Row 132
$ws.Range("H132").Value = 2182578.8
$ws.Range("I132").Value = 3581.3438
$ws.Range("K132").Value = 10744.0314
$ws.Range("M132").Value = -8214.0314

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 25280.352
$ws.Range("I20").Value = 7766.1333
$ws.Range("K20").Value = 7766.1333
$ws.Range("M20").Value = -7519.1333
# Row 105
$ws.Range("H105").Value = 5198.8
$ws.Range("I105").Value = 5198.8
$ws.Range("K105").Value = 5198.8
$ws.Range("M105").Value = -3451.8
# Row 113
$ws.Range("H113").Value = 23800.428
$ws.Range("I113").Value = 23800.428
$ws.Range("K113").Value = 23800.428
$ws.Range("M113").Value = -21630.428
# Row 134
$ws.Range("H134").Value = 9019.556
$ws.Range("I134").Value = 2244.8
$ws.Range("K134").Value = 6734.400000000001
$ws.Range("M134").Value = -4199.400000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
# Row 132
$ws.Range("H132").Value = 6997.875
$ws.Range("I132").Value = 2449.1052
$ws.Range("K132").Value = 7347.3156
$ws.Range("M132").Value = -4817.3156

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 2841882.5
$ws.Range("I107").Value = 803.3333
$ws.Range("J107").Value = 3907287.2
$ws.Range("K107").Value = 2409.9999
$ws.Range("L107").Value = 11721861.6
$ws.Range("M107").Value = -489.9998999999998
$ws.Range("N107").Value = -11725701.6
# Row 113
$ws.Range("H113").Value = 1326.4546
$ws.Range("J113").Value = 1156
$ws.Range("L113").Value = 3468
$ws.Range("N113").Value = -7808
# Row 114
$ws.Range("H114").Value = 1448.0625
$ws.Range("I114").Value = 1097
$ws.Range("J114").Value = 1721.1111
$ws.Range("K114").Value = 3291
$ws.Range("L114").Value = 5163.3333
$ws.Range("M114").Value = -37
$ws.Range("N114").Value = -11671.3333
# Row 121
$ws.Range("H121").Value = 1830
$ws.Range("I121").Value = 1000
$ws.Range("J121").Value = 1996
$ws.Range("K121").Value = 3000
$ws.Range("L121").Value = 5988
$ws.Range("M121").Value = -1690
$ws.Range("N121").Value = -8608
# Row 122
$ws.Range("H122").Value = 8611307
$ws.Range("I122").Value = 15572967
$ws.Range("J122").Value = 2185159.8
$ws.Range("K122").Value = 140156703
$ws.Range("L122").Value = 19666438.2
$ws.Range("M122").Value = -140154253
$ws.Range("N122").Value = -19671338.2
# Row 129
$ws.Range("H129").Value = 1925.4546
$ws.Range("J129").Value = 2516.1428
$ws.Range("L129").Value = 7548.428400000001
$ws.Range("N129").Value = -17548.4284
# Row 131
$ws.Range("H131").Value = 1493.26
$ws.Range("J131").Value = 1493.8469
$ws.Range("L131").Value = 4481.5407
$ws.Range("N131").Value = -14561.5407

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 9048.1
$ws.Range("I122").Value = 2926
$ws.Range("J122").Value = 23333
$ws.Range("K122").Value = 8778
$ws.Range("L122").Value = 69999
$ws.Range("M122").Value = -6328
$ws.Range("N122").Value = -74899
# Row 123
$ws.Range("H123").Value = 61635.445
$ws.Range("J123").Value = 61214.875
$ws.Range("L123").Value = 61214.875
$ws.Range("N123").Value = -66114.875
# Row 132
$ws.Range("H132").Value = 12520.315
$ws.Range("I132").Value = 7182.091
$ws.Range("J132").Value = 19860.375
$ws.Range("K132").Value = 21546.273
$ws.Range("L132").Value = 59581.125
$ws.Range("M132").Value = -19016.273
$ws.Range("N132").Value = -64641.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 12363.546
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888
# Row 16
$ws.Range("H16").Value = 993.1372699999999
$ws.Range("I16").Value = 924.53845
$ws.Range("K16").Value = 924.53845
$ws.Range("M16").Value = -754.53845
# Row 22
$ws.Range("H22").Value = 9569.25
$ws.Range("I22").Value = 7904.9414
$ws.Range("J22").Value = 13611.143
$ws.Range("K22").Value = 7904.9414
$ws.Range("L22").Value = 13611.143
$ws.Range("M22").Value = -7609.9414
$ws.Range("N22").Value = -14201.143
# Row 27
$ws.Range("H27").Value = 9569.25
$ws.Range("I27").Value = 7904.9414
$ws.Range("J27").Value = 13611.143
$ws.Range("K27").Value = 7904.9414
$ws.Range("L27").Value = 13611.143
$ws.Range("M27").Value = -7797.9414
$ws.Range("N27").Value = -13825.143
# Row 46
$ws.Range("H46").Value = 558787.6
$ws.Range("I46").Value = 771898.5600000001
$ws.Range("J46").Value = 4699.2
$ws.Range("K46").Value = 771898.5600000001
$ws.Range("L46").Value = 4699.2
$ws.Range("M46").Value = -771710.5600000001
$ws.Range("N46").Value = -5075.2
# Row 126
$ws.Range("H126").Value = 12363.546
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
# Row 133
$ws.Range("H133").Value = 79996
$ws.Range("J133").Value = 79996
$ws.Range("L133").Value = 79996
$ws.Range("N133").Value = -85056

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 4031.0527
$ws.Range("I107").Value = 1599.25
$ws.Range("J107").Value = 8199.857
$ws.Range("K107").Value = 4797.75
$ws.Range("L107").Value = 24599.571
$ws.Range("M107").Value = -2877.75
$ws.Range("N107").Value = -28439.571
# Row 122
$ws.Range("H122").Value = 83342536
$ws.Range("I122").Value = 333336160
$ws.Range("K122").Value = 1000008480
$ws.Range("M122").Value = -1000006030

